# Modifications to accomodate local builds
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 31 / column D: new, longer text replacing the old note.
$ws.Range("D31").Value = "Audio Amp and Codec Analog shut down when not in use"
$ws.Rows.Item(31).RowHeight = 28.8

# Row 75 and 77: status changes from "Open" to "Closed"
$ws.Range("C75").Value = "Closed"
$ws.Range("C77").Value = "Closed"

# New rows 80 and 81
$ws.Range("B80").Value = "State device name after successful setup"
$ws.Range("C80").Value = "Open"

$ws.Range("B81").Value = "Remove photon picture from device setup screen"
$ws.Range("C81").Value = "Open"

# Update the view to reflect the newly added rows
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 58
$ws.Range("C81").Select()
